$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.103.74"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "3.331.23"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.327.19"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("E11").Value = "  -1.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.83%  "

$ws.Range("E13").Value = "  -2.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "677.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.27%  "

$ws.Range("D15").Value = "3.861.33"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").Value = "66.285.65"
$ws.Range("E17").Value = "  -0.31%  "

# Rows 18 and 19 swap: Chainlink <-> TRON
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").Value = "3.327.03"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.34%  "

$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "

$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "608.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "

# Rows 35 and 36 swap: Hedera <-> Maker
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.834.20"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "

# Rows 40 and 41 swap: PEPE <-> Fetch.AI
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0702"
$ws.Range("E41").Value = "  -4.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "32.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.78%  "

$ws.Range("E45").Value = "  -2.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("E47").Value = "  -12.98%  "

$ws.Range("E48").Value = "  -2.10%  "

$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("E50").Value = "  -2.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.88%  "

